$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5104649250691508
$ws.Range("D2").Value = 0.7233556322385769
$ws.Range("E2").Value = 0.5057390294418187
$ws.Range("F2").Value = 0.1083173098324182
$ws.Range("G2").Value = -0.1163970892303022
$ws.Range("H2").Value = 0.1536126043573058
$ws.Range("I2").Value = -0.0251010488369374
$ws.Range("J2").Value = 0.1250703882501619
$ws.Range("K2").Value = 0.3996643215959416
$ws.Range("L2").Value = 0.5871125149530595
$ws.Range("M2").Value = 0.7747205868475892
$ws.Range("N2").Value = 0.7123937378756482
$ws.Range("O2").Value = 0.9360391759887446
$ws.Range("P2").Value = 0.06703388413220052
$ws.Range("D3").Value = 0.4735182507697232
$ws.Range("E3").Value = 0.7384320603883491
$ws.Range("F3").Value = 0.3832053538601286
$ws.Range("G3").Value = 0.05771416018395858
$ws.Range("H3").Value = 0.07536598196678382
$ws.Range("I3").Value = -0.2627283503193028
$ws.Range("J3").Value = -0.150248250641281
$ws.Range("K3").Value = -0.04201050655222459
$ws.Range("L3").Value = 0.0420347110983251
$ws.Range("M3").Value = 0.1932637093873051
$ws.Range("N3").Value = 0.05739492242263098
$ws.Range("O3").Value = 0.3068630469059259
$ws.Range("P3").Value = -0.1057120563392508
$ws.Range("E4").Value = 0.5815948407419296
$ws.Range("F4").Value = 0.1683111272338748
$ws.Range("G4").Value = 0.08075229237313496
$ws.Range("H4").Value = 0.106872662720764
$ws.Range("I4").Value = 0.139062134507934
$ws.Range("J4").Value = 0.2018692264075439
$ws.Range("K4").Value = 0.3446997744153577
$ws.Range("L4").Value = 0.4488833321568892
$ws.Range("M4").Value = 0.6989563963343225
$ws.Range("N4").Value = 0.4735459117008906
$ws.Range("O4").Value = 0.6579358914479844
$ws.Range("P4").Value = 0.08835831955228235
$ws.Range("F5").Value = 0.3166443726891301
$ws.Range("G5").Value = 0.03889375592201257
$ws.Range("H5").Value = -0.002106030724339357
$ws.Range("I5").Value = -0.2537524400723127
$ws.Range("J5").Value = -0.1189776463281428
$ws.Range("K5").Value = 0.005295257924840565
$ws.Range("L5").Value = 0.04497972503972599
$ws.Range("M5").Value = 0.2582749292304808
$ws.Range("N5").Value = 0.09736891141738298
$ws.Range("O5").Value = 0.3325753595824584
$ws.Range("P5").Value = -0.05726091530797562
$ws.Range("G6").Value = 0.4478117382092922
$ws.Range("H6").Value = 0.3260078053257223
$ws.Range("I6").Value = 0.1141288629224069
$ws.Range("J6").Value = 0.2058866507256234
$ws.Range("K6").Value = 0.2839052405011384
$ws.Range("L6").Value = 0.1468899967379731
$ws.Range("M6").Value = 0.08789781719038475
$ws.Range("N6").Value = -0.009311179259293425
$ws.Range("O6").Value = 0.007590821169478759
$ws.Range("P6").Value = 0.01686703050869204
$ws.Range("H7").Value = 0.3314966159157822
$ws.Range("I7").Value = 0.3085468373220643
$ws.Range("J7").Value = 0.3648456315955734
$ws.Range("K7").Value = 0.2972519644513208
$ws.Range("L7").Value = 0.05440784859433673
$ws.Range("M7").Value = 0.04901922246282474
$ws.Range("N7").Value = -0.08453128532741998
$ws.Range("O7").Value = -0.1322647584941255
$ws.Range("P7").Value = -0.06817278732926482
$ws.Range("I8").Value = 0.3457065951421486
$ws.Range("J8").Value = 0.4799415095463707
$ws.Range("K8").Value = 0.3029207710523865
$ws.Range("L8").Value = 0.2516860973554148
$ws.Range("M8").Value = 0.3811783426416966
$ws.Range("N8").Value = 0.1930620340375953
$ws.Range("O8").Value = 0.1692931569928255
$ws.Range("P8").Value = -0.06779520437147191
$ws.Range("J9").Value = 0.8029695880231121
$ws.Range("K9").Value = 0.5176673979838983
$ws.Range("L9").Value = 0.3894764285790648
$ws.Range("M9").Value = 0.3427623743421981
$ws.Range("N9").Value = 0.2195989319270962
$ws.Range("O9").Value = 0.09491040861790445
$ws.Range("P9").Value = 0.1237064202889188
$ws.Range("K10").Value = 0.684498687625048
$ws.Range("L10").Value = 0.5199013155728742
$ws.Range("M10").Value = 0.4702909441644824
$ws.Range("N10").Value = 0.395384492894999
$ws.Range("O10").Value = 0.2413549967551613
$ws.Range("P10").Value = 0.09452205310732144
$ws.Range("L11").Value = 0.7117584481777813
$ws.Range("M11").Value = 0.5956527406219603
$ws.Range("N11").Value = 0.6550540203292663
$ws.Range("O11").Value = 0.5209788817477098
$ws.Range("P11").Value = 0.0637430882359884
$ws.Range("M12").Value = 0.767226244536081
$ws.Range("N12").Value = 0.8079969154549849
$ws.Range("O12").Value = 0.7073950887642011
$ws.Range("P12").Value = 0.1252977586551207
$ws.Range("N13").Value = 0.7963264807952301
$ws.Range("O13").Value = 0.8439223777038533
$ws.Range("P13").Value = 0.1037459965656436
$ws.Range("O14").Value = 0.8431034520041651
$ws.Range("P14").Value = 0.1617483144037699
$ws.Range("P15").Value = 0.08518688770160042
